$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.07405227488138073
$ws.Range("C2").Value = 0.5849218257695413
$ws.Range("D2").Value = 0.5369781482230573
$ws.Range("E2").Value = 0.7327879285462182
$ws.Range("F2").Value = 0.7501725040260092
$ws.Range("G2").Value = 18

$ws.Range("B3").Value = -0.04979245339619409
$ws.Range("C3").Value = 0.6646959867222497
$ws.Range("D3").Value = 0.6330038130941004
$ws.Range("E3").Value = 0.7956153675577794
$ws.Range("F3").Value = 0.8158147246396246
$ws.Range("G3").Value = 19

$ws.Range("B4").Value = 0.2277465025845359
$ws.Range("C4").Value = 0.6805023410625829
$ws.Range("D4").Value = 0.6129072834672828
$ws.Range("E4").Value = 0.7828839527460522
$ws.Range("F4").Value = 0.770740616034928
$ws.Range("G4").Value = 18
